# Insert a new weekly record row above row 109 in the Espinaca price sheet.
# This pushes existing rows 109-154 down to 110-155, and populates the
# new row 109 with the latest weekly observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 109 (shifts rows 109:154 down to 110:155,
# Excel copies formatting/number-format from the row above automatically).
$ws.Rows.Item(109).Insert()

# Populate the newly inserted row 109 with its data.
$ws.Range("A109").Value = 8
$ws.Range("B109").Value = "Terminal La Palmera de La Serena"
$ws.Range("C109").Value = "Coquimbo"
$ws.Range("D109").Value = 44466
$ws.Range("E109").Value = 4
$ws.Range("F109").Value = 100112012
$ws.Range("G109").Value = "Espinaca"
$ws.Range("H109").Value = "Sin especificar"
$ws.Range("I109").Value = "Primera"
$ws.Range("J109").Value = 3320
$ws.Range("K109").Value = 400
$ws.Range("L109").Value = 500
$ws.Range("M109").Value = 450
$ws.Range("N109").Value = "`$/atado 300 a 500 gramos"
$ws.Range("O109").Value = "Provincia del Elquí"
$ws.Range("P109").Value = 900
$ws.Range("Q109").Value = 0.5
$ws.Range("R109").Value = "Hortaliza"
